# "merge translations for 4.2"
#
# The only functional change in this commit is the Russian footer string on
# the "Чат" sheet: the old text
#   "Документ сформирован $.Now, время в документе указано в часовом поясе $.Tz"
# is replaced by
#   "Сформировано $.BrandName в $.Now, время в документе указано в часовом поясе $.Tz"
# (i.e. the brand-name placeholder is added to the generated-on notice).
# The author also left the cursor/selection parked on that cell (A2) instead
# of B16, and nudged its font back to the sheet's standard body font
# (Calibri 12, top-vertical alignment) while retitling it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Чат")

$footerCell = $ws.Range("A2")
$footerCell.Value = "Сформировано `$.BrandName в `$.Now, время в документе указано в часовом поясе `$.Tz"

# Keep/confirm the body-text look (Calibri 12, top aligned) for the footer line.
$footerCell.Font.Name = "Calibri"
$footerCell.Font.Size = 12
$footerCell.VerticalAlignment = -4160  # xlVAlignTop

# Move the saved selection to the footer cell (was B16 before the edit).
$footerCell.Select()
